$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that receive updated vm_pu values for the 380 kV case (G, H, N unchanged)
$cols = @("B","C","D","E","F","I","J","K","L","M")

$rowData = @{
    2 = @(1.02, 1.044347406509429, 1.051968282271224, 1.052089468472426, 1.062542699770794, 1.042932735327767, 1.049412977925657, 1.054718315282301, 1.054839165913459, 1.065263763622047)
    3 = @(1.02, 1.045343851258234, 1.052749651937034, 1.052953231347063, 1.063439417389199, 1.043152968992481, 1.050056408429247, 1.055312517614185, 1.055515574008129, 1.065975111797181)
    4 = @(1.02, 1.045988991601384, 1.053255244276345, 1.053512762115652, 1.064020079858499, 1.043293857447578, 1.050472513907415, 1.055696341392798, 1.055953230913459, 1.066435190451662)
    5 = @(1.02, 1.046260297376032, 1.053467792565738, 1.05374813568778, 1.064264291108597, 1.043352699109338, 1.050647387016301, 1.055857540747686, 1.056137215074529, 1.066628555876498)
    6 = @(1.02, 1.046305855993346, 1.053503480161589, 1.053787664521422, 1.064305301122512, 1.043362556121698, 1.050676745575574, 1.055884597407832, 1.056168106397652, 1.066661019722323)
    7 = @(1.02, 1.045992616454539, 1.053258084369326, 1.053515906613406, 1.064023342626926, 1.043294645217492, 1.050474850798021, 1.055698495977046, 1.055955689346773, 1.066437774414216)
    8 = @(1.02, 1.044684081793976, 1.052232350303211, 1.052381252399295, 1.062845660666199, 1.043007498856112, 1.049630476840282, 1.054919265811352, 1.055067765658037, 1.065504210143885)
    9 = @(1.02, 1.042381173604827, 1.050424885170536, 1.050386636631546, 1.060773755391654, 1.042489149843759, 1.048140804726659, 1.053541114161126, 1.053502986677744, 1.063857580832341)
    10 = @(1.02, 1.040847891787883, 1.049219988837329, 1.049060184164096, 1.059394792829575, 1.042135309414605, 1.047146545637739, 1.052619009500178, 1.052459760146617, 1.062758835441816)
    11 = @(1.02, 1.040184444223569, 1.048698289728792, 1.048486610842391, 1.058798249518441, 1.041980137350508, 1.046715759332397, 1.052218948258388, 1.052008033421842, 1.062282844009538)
    12 = @(1.02, 1.039938082001035, 1.048504512981638, 1.048273679897675, 1.058576751582395, 1.041922206088461, 1.046555706998801, 1.052070230828976, 1.051840242471265, 1.06210600641574)
    13 = @(1.02, 1.039990924278011, 1.04854607848252, 1.048319348872115, 1.058624259784605, 1.041934645811293, 1.046590040525405, 1.052102136496814, 1.051876234164421, 1.062143940165663)
    14 = @(1.02, 1.040164078378432, 1.048682271962792, 1.048469007455163, 1.058779938692956, 1.041975354717443, 1.046702530150093, 1.052206657608276, 1.051994163754882, 1.062268227230085)
    15 = @(1.02, 1.040270773858907, 1.048766186029015, 1.048561232946199, 1.058875868880287, 1.042000397949518, 1.046771833560096, 1.052271041054696, 1.052066824156111, 1.062344800224944)
    16 = @(1.02, 1.040891932646057, 1.049254612999869, 1.049098267064022, 1.059434395272674, 1.04214556649175, 1.0471751299775, 1.052645543808797, 1.052489739778921, 1.062790420718127)
    17 = @(1.02, 1.041281696351277, 1.04956099892114, 1.049435346510898, 1.059784893976193, 1.042236103369134, 1.047428036617396, 1.052880250266552, 1.052755023461918, 1.063069886488516)
    18 = @(1.02, 1.041509084580057, 1.049739711360174, 1.049632035334215, 1.059989387561974, 1.042288723156613, 1.047575527106933, 1.053017074883754, 1.05290975862136, 1.063232872134748)
    19 = @(1.02, 1.041586625815507, 1.049800648108052, 1.04969911401084, 1.060059123649072, 1.042306633097041, 1.047625813173806, 1.053063715677082, 1.052962519230108, 1.063288442264966)
    20 = @(1.02, 1.041239873667849, 1.049528126322014, 1.049399173198278, 1.059747283227361, 1.04222640914782, 1.047400904778402, 1.052855076307222, 1.052726561060513, 1.063039904725557)
    21 = @(1.02, 1.040113086826385, 1.048642166215259, 1.048424933418493, 1.058734092783895, 1.041963375055491, 1.046669405832986, 1.05217588196515, 1.051959436409813, 1.062231628695473)
    22 = @(1.02, 1.039405046594322, 1.048085160511319, 1.04781308345524, 1.058097551015613, 1.041796297547179, 1.04620925746664, 1.051748169782242, 1.051477117838633, 1.061723241978731)
    23 = @(1.02, 1.039780352546443, 1.04838043629267, 1.0481373706067, 1.058434946886319, 1.041885029221599, 1.046453211979521, 1.051974971858642, 1.051732803447033, 1.061992765219268)
    24 = @(1.02, 1.041258771397167, 1.049542980026274, 1.049415518126184, 1.059764277740223, 1.042230790132662, 1.047413164567933, 1.052866451569757, 1.052739421995614, 1.063053452261753)
    25 = @(1.02, 1.04297618288188, 1.050892149468982, 1.050901717671391, 1.061308991305864, 1.042624616542085, 1.048526126073633, 1.053897991741944, 1.053907530809499, 1.064283453512311)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
